$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New "ParticipantsTab" query text (row 2, column B = "query").
# Replaces the old Cypher query with the updated version that adds
# diagnosis/genomic_info optional matches and sorts the collected samples.
$newParticipantsQuery = 'MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in [''TSV'']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'''') as `Participant ID`,
coalesce(s.study_name, '''') as `Study Name`,
coalesce(s.phs_accession,'''') as `Accession`,
coalesce(p.gender,'''') as `Gender`,
coalesce(apoc.text.join(samp, '',''), '''') as `Samples`
ORDER BY p.participant_id LIMIT 100'

$ws.Range("B2").Value = $newParticipantsQuery
